$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet / workbook tab
$ws.Name = "Controle de temperatura"

# Start from a clean sheet so no stale cells from the old layout remain
$ws.Cells.Clear()

# Title row
$ws.Range("A1").Value = "Farmácia Rio Negro - Silva & Heidrich ltda"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Size = 16

# Row 2 intentionally left blank (spacer)

# Subtitle row
$ws.Range("A3").Value = "Controle de temperatura - Sensor DHT-22"

# Row 4 intentionally left blank (spacer)

# Header row
$ws.Range("A5").Value = "Data"
$ws.Range("B5").Value = "Temperatura (ºC)"
$ws.Range("C5").Value = "Umidade (%)"
$ws.Range("D5").Value = "Hora da Leitura"
$ws.Range("E5").Value = "Temperatura Média (ºC)"
$ws.Range("A5:E5").Font.Bold = $true

# Data row - force column A to be stored as plain text so the
# dd/mm/yyyy reading isn't reinterpreted as a date serial number
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "04/01/2023"
$ws.Range("A6").Style = "Normal"

$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = "20:56:00"
$ws.Range("E6").Value = ""

# Column widths (values chosen so the saved file reports the exact
# target character widths after Excel's pixel-rounding conversion)
$ws.Columns.Item(1).ColumnWidth = 38.165
$ws.Columns.Item(2).ColumnWidth = 15.165
$ws.Columns.Item(4).ColumnWidth = 14.165
$ws.Columns.Item(5).ColumnWidth = 19.165
